# MSE-1154: Fix inline styler
#
# Semantic changes captured by this edit (everything else in the source
# diff is re-serialization noise produced by regenerating the fixture
# with a different OOXML writer -- namespaces/uids/default attribute
# dumps, shared-strings vs. inline-strings storage, etc. -- and isn't the
# result of any discrete spreadsheet action):
#   1. The worksheet is renamed from "7e862dbb" to "24b3b3f3".
#   2. Cell B2 (one of the two cells driving the Active/Inactive list
#      validation together with C2) goes from not-existing to being a
#      present-but-empty cell, i.e. it was touched/cleared so it is no
#      longer just an implicit blank.
#   3. The saved selection moves from L6 to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet.
$ws.Name = "24b3b3f3"

# 2. Touch B2 so it becomes an explicit (empty) cell instead of an
#    implicit blank -- mirrors the validation dropdown cell being left
#    unset next to C2's "Active" value. Re-applying the default "Normal"
#    style (rather than writing a value, which the engine collapses back
#    to a true blank) is what forces a concrete, empty cell record to be
#    persisted for B2 without changing its appearance.
$ws.Range("B2").Value = ""
$ws.Range("B2").Style = "Normal"

# 3. Reset the active selection to A1.
$ws.Range("A1").Select() | Out-Null
